# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before "总计", populated with
#    the per-fund holdings snapshot for 2022-Q1 (same layout as the other
#    quarterly sheets).
# 2. Insert a new row at the top of the "总计" (totals) sheet's data with
#    the 2022-Q1 aggregate row, shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing text storage (so values
# that look numeric, like "96.71" or "001643", stay strings instead of
# being auto-coerced to numbers).
# ---------------------------------------------------------------------
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

# =======================================================================
# Step 1: add the "2022-Q1" sheet before "总计"
# =======================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Copy header-row formatting from an existing quarterly sheet so the new
# sheet's header cells (B1:H1) pick up the same bold/border/center style.
$headerTemplate = $wb.Worksheets.Item("2021-Q4")
$headerTemplate.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Copy the "index" column (A) style from an existing sheet's data rows so
# A2:A9 match the A2-style (bold/centered/bordered) used elsewhere.
$headerTemplate.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

$q1Data = @(
    @(0, "540008", "汇丰晋信低碳先锋股票",   "96.71", "93.08", "3.72", "3.5976", 9),
    @(1, "001643", "汇丰晋信智造先锋股票A",  "29.09", "92.99", "3.50", "1.0182", 10),
    @(2, "007994", "华夏中证500指数增强A",   "31.45", "92.72", "1.32", "0.4151", 10),
    @(3, "001644", "汇丰晋信智造先锋股票C",  "10.91", "92.99", "3.50", "0.3818", 10),
    @(4, "007995", "华夏中证500指数增强C",   "5.45",  "92.72", "1.32", "0.0719", 10),
    @(5, "013204", "恒生前海恒源天利债A",    "1.29",  "21.66", "1.08", "0.0139", 9),
    @(6, "001797", "华融新利灵活配置混合",   "0.02",  "48.66", "2.70", "0.0005", 4),
    @(7, "013205", "恒生前海恒源天利债C",    "0.00",  "21.66", "1.08", 0,        9)
)

$row = 2
foreach ($rec in $q1Data) {
    $q1.Range("A$row").Value = $rec[0]
    Set-TextCell $q1 "B$row" $rec[1]
    Set-TextCell $q1 "C$row" $rec[2]
    Set-TextCell $q1 "D$row" $rec[3]
    Set-TextCell $q1 "E$row" $rec[4]
    Set-TextCell $q1 "F$row" $rec[5]
    if ($row -eq 9) {
        $q1.Range("G$row").Value = $rec[6]
    } else {
        Set-TextCell $q1 "G$row" $rec[6]
    }
    $q1.Range("H$row").Value = $rec[7]
    $row = $row + 1
}

# =======================================================================
# Step 2: insert the 2022-Q1 summary row at the top of "总计"'s data
# =======================================================================
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

# Restore the "index" column style (A2) to match the other data rows, and
# clear the stray formatting the row-insert applied to B2:D2.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2:D2").ClearFormats()

$ws.Range("A2").Value = 0
Set-TextCell $ws "B2" "2022-Q1"
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 5.5

# Renumber the index column for the rows that shifted down
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
